$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.847.60'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.42%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.336.13'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -3.28%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.30%  '
$ws.Range('E7').Value = '  -3.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.326.79'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.37%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.611'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.77%  '
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '53.65'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.42%  '
$ws.Range('E13').Value = '  -2.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.87'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.00%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.784.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.27%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.13'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.332.37'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.08%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.33%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '63.798.75'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.975'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.53%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '409.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.90%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.03'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.38'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +6.00%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '13.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +10.24%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.96'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.51%  '
$ws.Range('E27').Value = '  -2.86%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.73'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.60'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -5.77%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.04'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -3.26%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.37'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.91%  '
$ws.Range('B32').Value = 'Bittensor'
$ws.Range('C32').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '579.35'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.53%  '
$ws.Range('B33').Value = 'Cosmos'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.30%  '
$ws.Range('E34').Value = '  -3.03%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '57.84'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.11%  '
$ws.Range('E36').Value = '  +1.79%  '
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.17'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -6.16%  '
$ws.Range('E39').Value = '  +2.29%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0₃0738'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.25%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.366'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -3.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.131.81'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.83%  '
$ws.Range('E43').Value = '  -0.19%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.80'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.26'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0400'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.43%  '
$ws.Range('E47').Value = '  -5.52%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -3.93%  '
$ws.Range('E49').Value = '  -3.82%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '132.25'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -4.36%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.06'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.32%  '
